# Generate Report for Archive
#
# The localization status for both language targets moved from
# "Ready for handoff" to "In Translation": update every cell that shows the
# status (the Overview rollup columns zh-cn/de-de, plus the Status column on
# each per-locale detail sheet), then shrink those now-narrower columns to
# match the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn status in column E, de-de status in column F (rows 2-3)
$overview.Range("E2:F3").Value = $newStatus

# Per-locale detail sheets: Status column is column C (rows 2-3)
$zhcn.Range("C2:C3").Value = $newStatus
$dede.Range("C2:C3").Value = $newStatus

# The status text got shorter, so the columns that hold it get narrower too.
$overview.Columns.Item(5).ColumnWidth = 12.5   # E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # C (Status)
